$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add rows 195-203 below the existing data, reusing the formatting of the
#     rows immediately above (row 192 -> date column style 10 / value columns
#     style 11; row 193 -> date column style 10 / value columns style 1).
#     Done before the P193 style tweak below so the copied formatting is the
#     original (unmodified) one. ---
$ws.Range("A192:Q192").Copy()
$ws.Range("A195:Q202").PasteSpecial(-4122)
$ws.Range("A193:Q193").Copy()
$ws.Range("A203:Q203").PasteSpecial(-4122)

# --- Change P193's formatting to a distinct style (style 13 in the target) ---
$ws.Range("P193").Font.Name = $null

$rowValues = @(44145,925918,139455,785426,1037,2250,133360,241,32,48,0,1,851,1195,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(195, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44146,929945,139727,789746,472,2257,133846,225,26,46,0,1,857,1196,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(196, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44147,933933,139954,793449,530,2258,134170,219,31,47,0,2,857,1197,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(197, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44148,937575,140175,796520,880,2260,134213,216,30,46,0,2,857,1199,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(198, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44149,941967,140474,800652,841,2263,134480,213,30,49,0,2,858,1201,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(199, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44150,945113,140635,803756,722,2266,134691,214,27,49,0,2,859,1203,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(200, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44151,946579,140740,805119,720,2268,134905,201,28,48,0,3,860,1204,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(201, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44152,950955,140953,808279,1723,2272,135233,205,34,46,0,3,862,1206,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(202, $col).Value = $rowValues[$col - 1] }

$rowValues = @(44153,954878,141196,813159,523,2278,135511,194,16,44,0,3,866,1208,39,123,42)
for ($col = 1; $col -le 17; $col++) { $ws.Cells.Item(203, $col).Value = $rowValues[$col - 1] }
